$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 - "Save" - formatted like the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New data column H2:H13 - 0/1 "Save" flags
$saveValues = @(0, 0, 1, 0, 1, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
